# Update "想去人数" (want-to-go count) values on the 展览 (Exhibition),
# 演出 (Performance) and 全部类型 (All types) sheets to reflect newly
# generated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 8627
$wsExpo.Range("F7").Value = 10857
$wsExpo.Range("F24").Value = 583
$wsExpo.Range("F32").Value = 5
$wsExpo.Range("F37").Value = 347

# --- Sheet "演出" ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F17").Value = 390

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 8627
$wsAll.Range("F11").Value = 10857
$wsAll.Range("F21").Value = 583
$wsAll.Range("F37").Value = 347
$wsAll.Range("F45").Value = 390
